$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.595.24"
$ws.Range("E2").Value = "  -1.45%  "
$ws.Range("D3").Value = "'1.751.26"
$ws.Range("E3").Value = "  -0.74%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'324.04"
$ws.Range("E5").Value = "  +1.03%  "
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").Value = "'0.4597"
$ws.Range("E7").Value = "  +8.12%  "
$ws.Range("D8").Value = "'0.3592"
$ws.Range("E8").Value = "  -0.61%  "
$ws.Range("D9").Value = "'0.07500"
$ws.Range("E9").Value = "  +0.70%  "
$ws.Range("E10").Value = "  -3.80%  "
$ws.Range("D11").Value = "'1.096"
$ws.Range("E11").Value = "  -0.34%  "
$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("D13").Value = "'20.72"
$ws.Range("E13").Value = "  -2.28%  "
$ws.Range("D14").Value = "'6.005"
$ws.Range("E14").Value = "  -1.31%  "
$ws.Range("D15").Value = "'7.093"
$ws.Range("E15").Value = "  -3.14%  "
$ws.Range("D16").Value = "'1.756.84"
$ws.Range("E16").Value = "  -0.61%  "
$ws.Range("D17").Value = "'92.32"
$ws.Range("E17").Value = "  +1.44%  "
$ws.Range("D18").Value = "'0.00001066"
$ws.Range("E18").Value = "  +0.82%  "
$ws.Range("D19").Value = "'0.06405"
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("D21").Value = "'16.74"
$ws.Range("E21").Value = "  -1.71%  "
$ws.Range("D22").Value = "'5.822"
$ws.Range("D23").Value = "'27.638.34"
$ws.Range("E23").Value = "  -1.37%  "
$ws.Range("D24").Value = "'11.20"
$ws.Range("E24").Value = "  -0.81%  "
$ws.Range("D25").Value = "'2.114"
$ws.Range("E25").Value = "  -1.11%  "
$ws.Range("D26").Value = "'163.96"
$ws.Range("E26").Value = "  +3.83%  "
$ws.Range("D27").Value = "'20.41"
$ws.Range("E27").Value = "  +1.20%  "
$ws.Range("D28").Value = "'1.956.68"
$ws.Range("E28").Value = "  -1.43%  "
$ws.Range("D29").Value = "'2.088"
$ws.Range("E29").Value = "  -2.19%  "
$ws.Range("D30").Value = "'126.66"
$ws.Range("E30").Value = "  +1.69%  "
$ws.Range("D31").Value = "'1.074"
$ws.Range("E31").Value = "  -7.84%  "
$ws.Range("D32").Value = "'0.09218"
$ws.Range("E32").Value = "  +3.77%  "
$ws.Range("D33").Value = "'3.674"
$ws.Range("E33").Value = "  +3.51%  "
$ws.Range("D34").Value = "'5.523"
$ws.Range("E34").Value = "  -2.13%  "
$ws.Range("D35").Value = "'11.91"
$ws.Range("E35").Value = "  -5.09%  "
$ws.Range("D36").Value = "'0.02295"
$ws.Range("E36").Value = "  -1.03%  "
$ws.Range("D37").Value = "'0.2102"
$ws.Range("E37").Value = "  -0.36%  "
$ws.Range("E38").Value = "  -0.44%  "
$ws.Range("D39").Value = "'0.6360"
$ws.Range("E39").Value = "  -0.27%  "
$ws.Range("D40").Value = "'4.960"
$ws.Range("E40").Value = "  -1.29%  "
$ws.Range("D41").Value = "'1.200"
$ws.Range("E41").Value = "  +1.22%  "
$ws.Range("D42").Value = "'1.384"
$ws.Range("E42").Value = "  -0.85%  "
$ws.Range("D43").Value = "'7.787"
$ws.Range("E43").Value = "  -0.72%  "
$ws.Range("D44").Value = "'13.22"
$ws.Range("E44").Value = "  -2.26%  "
$ws.Range("D45").Value = "'0.5898"
$ws.Range("E45").Value = "  -0.30%  "
$ws.Range("D46").Value = "'3.707"
$ws.Range("E46").Value = "  +0.49%  "
$ws.Range("D47").Value = "'122.99"
$ws.Range("E47").Value = "  +0.19%  "
$ws.Range("D48").Value = "'1.952"
$ws.Range("E48").Value = "  -2.61%  "
$ws.Range("D49").Value = "'1.144"
$ws.Range("E49").Value = "  -3.70%  "
$ws.Range("D50").Value = "'0.06850"
$ws.Range("E50").Value = "  -0.26%  "
$ws.Range("D51").Value = "'72.12"
$ws.Range("E51").Value = "  -2.77%  "
